# Updated High CCS Storage Cost file, per proposal 129
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("high_ccs_storage_cost")

# Update period values for existing rows and append new rows (6-23)
$ws.Range("A6").Value = "ALL"
$ws.Range("B6").Value = "carbon-storage"
$ws.Range("C6").Value = "carbon-storage regional"
$ws.Range("D6").Value = "carbon-storage regional"
$ws.Range("E6").Value = 1975
$ws.Range("F6").Value = 10000
$ws.Range("A7").Value = "ALL"
$ws.Range("B7").Value = "carbon-storage"
$ws.Range("C7").Value = "carbon-storage regional"
$ws.Range("D7").Value = "carbon-storage regional"
$ws.Range("E7").Value = 1990
$ws.Range("F7").Value = 10000
$ws.Range("A8").Value = "ALL"
$ws.Range("B8").Value = "carbon-storage"
$ws.Range("C8").Value = "carbon-storage regional"
$ws.Range("D8").Value = "carbon-storage regional"
$ws.Range("E8").Value = 2005
$ws.Range("F8").Value = 10000
$ws.Range("A9").Value = "ALL"
$ws.Range("B9").Value = "carbon-storage"
$ws.Range("C9").Value = "carbon-storage regional"
$ws.Range("D9").Value = "carbon-storage regional"
$ws.Range("E9").Value = 2020
$ws.Range("F9").Value = 10000
$ws.Range("A10").Value = "ALL"
$ws.Range("B10").Value = "carbon-storage"
$ws.Range("C10").Value = "carbon-storage regional"
$ws.Range("D10").Value = "carbon-storage regional"
$ws.Range("E10").Value = 2035
$ws.Range("F10").Value = 10000
$ws.Range("A11").Value = "ALL"
$ws.Range("B11").Value = "carbon-storage"
$ws.Range("C11").Value = "carbon-storage regional"
$ws.Range("D11").Value = "carbon-storage regional"
$ws.Range("E11").Value = 2050
$ws.Range("F11").Value = 10000
$ws.Range("A12").Value = "ALL"
$ws.Range("B12").Value = "carbon-storage"
$ws.Range("C12").Value = "carbon-storage regional"
$ws.Range("D12").Value = "carbon-storage regional"
$ws.Range("E12").Value = 2065
$ws.Range("F12").Value = 10000
$ws.Range("A13").Value = "ALL"
$ws.Range("B13").Value = "carbon-storage"
$ws.Range("C13").Value = "carbon-storage regional"
$ws.Range("D13").Value = "carbon-storage regional"
$ws.Range("E13").Value = 2080
$ws.Range("F13").Value = 10000
$ws.Range("A14").Value = "ALL"
$ws.Range("B14").Value = "carbon-storage"
$ws.Range("C14").Value = "carbon-storage regional"
$ws.Range("D14").Value = "carbon-storage regional"
$ws.Range("E14").Value = 2095
$ws.Range("F14").Value = 10000
$ws.Range("A15").Value = "ALL"
$ws.Range("B15").Value = "carbon-storage"
$ws.Range("C15").Value = "carbon-storage regional"
$ws.Range("D15").Value = "carbon-storage regional"
$ws.Range("E15").Value = 1975
$ws.Range("F15").Value = 10000
$ws.Range("A16").Value = "ALL"
$ws.Range("B16").Value = "carbon-storage"
$ws.Range("C16").Value = "carbon-storage regional"
$ws.Range("D16").Value = "carbon-storage regional"
$ws.Range("E16").Value = 1990
$ws.Range("F16").Value = 10000
$ws.Range("A17").Value = "ALL"
$ws.Range("B17").Value = "carbon-storage"
$ws.Range("C17").Value = "carbon-storage regional"
$ws.Range("D17").Value = "carbon-storage regional"
$ws.Range("E17").Value = 2005
$ws.Range("F17").Value = 10000
$ws.Range("A18").Value = "ALL"
$ws.Range("B18").Value = "carbon-storage"
$ws.Range("C18").Value = "offshore carbon-storage"
$ws.Range("D18").Value = "offshore carbon-storage"
$ws.Range("E18").Value = 2020
$ws.Range("F18").Value = 10000
$ws.Range("A19").Value = "ALL"
$ws.Range("B19").Value = "carbon-storage"
$ws.Range("C19").Value = "offshore carbon-storage"
$ws.Range("D19").Value = "offshore carbon-storage"
$ws.Range("E19").Value = 2035
$ws.Range("F19").Value = 10000
$ws.Range("A20").Value = "ALL"
$ws.Range("B20").Value = "carbon-storage"
$ws.Range("C20").Value = "offshore carbon-storage"
$ws.Range("D20").Value = "offshore carbon-storage"
$ws.Range("E20").Value = 2050
$ws.Range("F20").Value = 10000
$ws.Range("A21").Value = "ALL"
$ws.Range("B21").Value = "carbon-storage"
$ws.Range("C21").Value = "offshore carbon-storage"
$ws.Range("D21").Value = "offshore carbon-storage"
$ws.Range("E21").Value = 2065
$ws.Range("F21").Value = 10000
$ws.Range("A22").Value = "ALL"
$ws.Range("B22").Value = "carbon-storage"
$ws.Range("C22").Value = "offshore carbon-storage"
$ws.Range("D22").Value = "offshore carbon-storage"
$ws.Range("E22").Value = 2080
$ws.Range("F22").Value = 10000
$ws.Range("A23").Value = "ALL"
$ws.Range("B23").Value = "carbon-storage"
$ws.Range("C23").Value = "offshore carbon-storage"
$ws.Range("D23").Value = "offshore carbon-storage"
$ws.Range("E23").Value = 2095
$ws.Range("F23").Value = 10000

# Make this sheet the active/selected tab and set the new selection
$ws.Activate()
$ws.Range("E18").Select()

# Match page setup (portrait orientation) applied in the source edit
$ws.PageSetup.Orientation = [Microsoft.Office.Interop.Excel.XlPageOrientation]::xlPortrait
